# Add a new "2022-Q4" worksheet right after "总计" and push every other
# quarterly sheet one position to the right, then update the "总计"
# summary sheet with a new row for the 2022-Q4 totals.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new detail sheet "2022-Q4" and position it right after
#    "总计" (position 1), before "2022-Q3".
# ---------------------------------------------------------------------
$zongji = $wb.Worksheets.Item(1)

$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "2022-Q4"
$newSheet.Move($null, $zongji)

# The runtime's worksheet handle can go stale across a Move(), so re-fetch
# it by name before writing any cell data into it.
$ws = $wb.Worksheets.Item("2022-Q4")

# Header row (bold/bordered style copied from the equivalent header on the
# "2020-Q4" sheet).
$ws.Cells.Item(1, 2).Value = "基金代码"
$ws.Cells.Item(1, 3).Value = "基金名称"
$ws.Cells.Item(1, 4).Value = "基金规模"
$ws.Cells.Item(1, 5).Value = "股票总仓位"
$ws.Cells.Item(1, 6).Value = "仓位占比"
$ws.Cells.Item(1, 7).Value = "持有市值(亿元)"
$ws.Cells.Item(1, 8).Value = "仓位排名"

$rows = @(
  @("217005", "招商先锋混合",               "8.79",  "72.35", "3.11", "0.2734", 10),
  @("010418", "财通景气行业混合A",           "2.60",  "86.37", "5.87", "0.1526", 6),
  @("501015", "财通多策略升级混合（LOF）A",   "1.97",  "86.66", "5.84", "0.1150", 7),
  @("005959", "财通新视野灵活配置混合C",      "1.61",  "86.51", "5.64", "0.0908", 7),
  @("005851", "财通新视野灵活配置混合A",      "0.63",  "86.51", "5.64", "0.0355", 7),
  @("015271", "财通多策略升级混合（LOF）C",   "0.57",  "86.66", "5.84", "0.0333", 7),
  @("016234", "财通景气行业混合C",           "0.00",  "86.37", "5.87", 0,        6)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
  $r = $i + 2
  $entry = $rows[$i]
  $ws.Cells.Item($r, 1).Value = $i
  $ws.Cells.Item($r, 2).Value = $entry[0]
  $ws.Cells.Item($r, 3).Value = $entry[1]
  $ws.Cells.Item($r, 4).Value = $entry[2]
  $ws.Cells.Item($r, 5).Value = $entry[3]
  $ws.Cells.Item($r, 6).Value = $entry[4]
  $ws.Cells.Item($r, 7).Value = $entry[5]
  $ws.Cells.Item($r, 8).Value = $entry[6]
}

# Copy the header style (bold font + border, "s=1" on the sibling sheets)
# and the column-A index style ("s=1") from the "2020-Q4" sheet so the new
# sheet matches the look of its peers.
$template = $wb.Worksheets.Item("2020-Q4")
$template.Range("B1:H1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)
$template.Range("A2").Copy()
$ws.Range("A2:A8").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2. Update the "总计" overview sheet: insert the 2022-Q4 totals as the
#    new row 2 and shift the rest of the quarters down by one row.
# ---------------------------------------------------------------------
$summary = @(
  @(0, "2022-Q4", 7,  0.7),
  @(1, "2022-Q3", 32, 6.64),
  @(2, "2022-Q2", 38, 9.109999999999999),
  @(3, "2022-Q1", 92, 22.58),
  @(4, "2021-Q4", 29, 6.55),
  @(5, "2021-Q3", 5,  0.28),
  @(6, "2021-Q2", 1,  0.04),
  @(7, "2021-Q1", 5,  0.46),
  @(8, "2020-Q4", 8,  0.78)
)

for ($i = 0; $i -lt $summary.Length; $i++) {
  $r = $i + 2
  $entry = $summary[$i]
  $zongji.Cells.Item($r, 1).Value = $entry[0]
  $zongji.Cells.Item($r, 2).Value = $entry[1]
  $zongji.Cells.Item($r, 3).Value = $entry[2]
  $zongji.Cells.Item($r, 4).Value = $entry[3]
}

# Row 10 is brand new (the sheet used to stop at row 9), so its column-A
# cell needs the same style ("s=2") as every other index cell above it.
$zongji.Cells.Item(9, 1).Copy()
$zongji.Cells.Item(10, 1).PasteSpecial(-4122)
$zongji.Cells.Item(10, 1).Value = 8
